$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a brand new row above the old row 7 ("拿走/nv/-"), which shifts the
# old rows 7, 8, 9 down to 8, 9, 10.
# ---------------------------------------------------------------------------
$ws.Rows("7:7").Insert()
$ws.Rows("7:7").RowHeight = 16.5

# ---------------------------------------------------------------------------
# Populate the freshly inserted row 7 with the new vocabulary entry.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "被"
$ws.Range("B7").Value = "passive"
$ws.Range("C7").Value = "虛詞"

# ---------------------------------------------------------------------------
# A7 should look exactly like the red/left-top-bottom-bordered style that is
# already used by column B of the "拿走/nv/-" row (now row 8). Re-using that
# formatting via copy/paste keeps the underlying style table clean.
# ---------------------------------------------------------------------------
$ws.Range("B8").Copy()
$ws.Range("A7").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------------
# C7 re-uses the right/top/bottom red border that column C of the same row
# (now row 8) already has, then gets recoloured to the new red font.
# ---------------------------------------------------------------------------
$ws.Range("C8").Copy()
$ws.Range("C7").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = 0

$ws.Range("C7").Font.Color = 255
$ws.Range("C7").Font.Name = "新細明體"
$ws.Range("C7").Font.Size = 11

# ---------------------------------------------------------------------------
# B7 gets the new red font plus a medium red border on just the top and
# bottom edges (no left/right border).
# ---------------------------------------------------------------------------
$ws.Range("B7").Font.Color = 255
$ws.Range("B7").Font.Name = "新細明體"
$ws.Range("B7").Font.Size = 11

$ws.Range("B7").Borders.Item(8).Weight = -4138
$ws.Range("B7").Borders.Item(8).Color = 255
$ws.Range("B7").Borders.Item(9).Weight = -4138
$ws.Range("B7").Borders.Item(9).Color = 255

# ---------------------------------------------------------------------------
# Reflect the same selection state Excel would leave behind after this edit.
# ---------------------------------------------------------------------------
$ws.Range("A7:C7").Select() | Out-Null
